$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 18 ("Linguagens de Programação" - Text Placeholder 2): rewrite body
# text to sz=2400, add a bold "Vídeo" lead-in, move the existing targethd.net
# hyperlink paragraph up, and append two new bullet paragraphs (Alan Turing /
# John Von Neumann) each with a bold word and a new YouTube hyperlink.
# ---------------------------------------------------------------------------
$s18 = $p.Slides.Item(18)
$shape = $s18.Shapes.Item(2)
$tf = $shape.TextFrame
$tr = $tf.TextRange

$fullText = "- Vídeo`r`thttps://www.targethd.net/a-evolucao-das-linguagens-de-programacao-entre-1965-e-2019/`r`r- Alan Turing – Pai da Computação (Software)`r`thttps://www.youtube.com/watch?v=5jAq6yU8bxg`r- John Von Neumann – Arquitetura de Computador (Hardware)"

$tr.Text = $fullText

# Apply the common font size to the whole text range first (keeps the
# existing Times New Roman / solidFill / lang inherited from the original
# run formatting).
$tr.Font.Size = 24

# "Vídeo" -> bold
$tr.Characters(3, 5).Font.Bold = -1

# targethd.net hyperlink (reuses existing rId3 relationship - same address)
$tr.Characters(10, 84).ActionSettings(1).Hyperlink.Address = "https://www.targethd.net/a-evolucao-das-linguagens-de-programacao-entre-1965-e-2019/"

# "Software" -> bold
$tr.Characters(131, 8).Font.Bold = -1

# YouTube hyperlink (new relationship, rId4)
$tr.Characters(142, 43).ActionSettings(1).Hyperlink.Address = "https://www.youtube.com/watch?v=5jAq6yU8bxg"

# "Hardware" -> bold
$tr.Characters(234, 8).Font.Bold = -1

# Re-touch the lone blank paragraph (between the two hyperlink paragraphs)
# so its run carries the same full rPr (lang/dirty/solidFill/latin/cs) as
# the rest of the body instead of only "sz".
$tr.Paragraphs(3).Text = ""

# ---------------------------------------------------------------------------
# Slide 2 title: "Aula 03" -> "Aula "
# ---------------------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$titleShape = $s2.Shapes.Item(4)
$titleTr = $titleShape.TextFrame.TextRange
$titleTr.Replace("Aula 03", "Aula ") | Out-Null
